$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.825.90"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "2.264.85"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.529"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.82%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.68%  "

$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.112"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("D14").Value = "2.615.63"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("D16").Value = "2.272.79"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.784"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.95%  "

$ws.Range("D18").Value = "41.757.05"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.62%  "

$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.04%  "

$ws.Range("E29").Value = "  -5.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.13%  "

$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0743"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.30%  "

$ws.Range("E37").Value = "  +1.56%  "

$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("E40").Value = "  +0.49%  "

$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").Value = "2.009.20"
$ws.Range("E43").Value = "  -2.30%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.39%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.59%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.62%  "

$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.32%  "

$ws.Range("E50").Value = "  +0.86%  "

$ws.Range("E51").Value = "  -0.07%  "
